# v3.0.0 : work from A to Z --> TODO : next steps is to then Read the Assembly list option
#
# The "Options" sheet's reference list is rewritten: a couple of labels are
# corrected/renamed, several obsolete references are dropped, and the
# "J900"/"J913"/"S643" rows are moved back into their proper sorted position.
# The stray time-format style that had been applied to B12 is also cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options")

# Final, corrected list of references for column A, starting at row 2
# (row 1 is the header: Référence / FROM / SRAM / DRAM).
$values = @(
    "J858",
    "J900",
    "H580#40A0",
    "H590#90E0",
    "0303~J535#655B",
    "J541#60V8",
    "J562#6552",
    "J571#656E",
    "J572#656F",
    "J5994#6569",
    "H501#G103",
    "H510#B11",
    "H990#32K",
    "J674",
    "J721",
    "J728",
    "J733",
    "J734",
    "J736",
    "J738#2M",
    "J749",
    "J801",
    "J802#11",
    "J803",
    "J804",
    "J807",
    "J819",
    "J829",
    "J830",
    "J835",
    "J836",
    "J842",
    "J848",
    "J850",
    "J853",
    "J854",
    "J873",
    "J876",
    "J890",
    "J893",
    "J894",
    "J895",
    "J900",
    "J913",
    "J917",
    "J930",
    "J937",
    "J948",
    "J953",
    "J956",
    "J971",
    "J965",
    "J981",
    "R955",
    "S617",
    "S643",
    "S656",
    "S661",
    "S707",
    "S728",
    "S731",
    "S837#1",
    "S985"
)

$firstDataRow = 2
$lastDataRowBefore = $ws.UsedRange.Rows.Count
$newLastDataRow = $firstDataRow + $values.Count - 1

# Write the corrected values into column A.
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($firstDataRow + $i, 1).Value = $values[$i]
}

# The previous list had more rows than the corrected one (several stale
# references were removed outright) -- delete the now-unused trailing rows.
if ($lastDataRowBefore -gt $newLastDataRow) {
    $deleteRange = $ws.Range($ws.Cells.Item($newLastDataRow + 1, 1), $ws.Cells.Item($lastDataRowBefore, 1))
    $deleteRange.EntireRow.Delete()
}

# B12 had a stray time number-format (h:mm) applied with no real value in it;
# clear that leftover cell entirely (value + formatting).
$ws.Cells.Item(12, 2).Clear()

# Force a full recalculation the next time the workbook is opened.
$wb.ForceFullCalculation = $true
